$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2BDS"
$ws.Range("B7").Value = 1203945
